# feat: add 2022-Q4 data
#
# - Inserts a new "2022-Q4" detail sheet (fund holdings) positioned
#   between "总计" and "2022-Q3".
# - Updates the "总计" (summary) sheet: a new row for 2022-Q4 is added at
#   row 2 (pushing the existing 2022-Q3 summary row down to row 3).

# Helper: force a value to be written as TEXT (not auto-coerced to a
# number), e.g. fund codes like "005269" or decimal-looking strings like
# "3.38" that must stay literal text. Excel's text-coercion trick (leading
# apostrophe) sets a "quote prefix" style bit, so we immediately reset the
# cell style back to Normal to avoid leaving stray formatting behind.
function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$sheetTotal = $wb.Worksheets.Item("总计")

### 1. Update "总计": insert the 2022-Q4 row, push 2022-Q3 row down to row 3 ###

# Move the existing row 2 (2022-Q3) data down to row 3 first.
Set-TextValue $sheetTotal.Cells.Item(3, 2) "2022-Q3"
$sheetTotal.Cells.Item(3, 3).Value = 1
$sheetTotal.Cells.Item(3, 4).Value = 0.01

# Row 3 col A needs the same styling as col A elsewhere in the column.
$sheetTotal.Cells.Item(2, 1).Copy()
$sheetTotal.Cells.Item(3, 1).PasteSpecial(-4122)  # xlPasteFormats
$sheetTotal.Cells.Item(3, 1).Value = 1

# Now write the new 2022-Q4 summary row into row 2.
Set-TextValue $sheetTotal.Cells.Item(2, 2) "2022-Q4"
$sheetTotal.Cells.Item(2, 3).Value = 3
$sheetTotal.Cells.Item(2, 4).Value = 0.28
# (A2 already holds 0 with the correct style - unchanged.)

### 2. Add the new "2022-Q4" sheet right after "总计" (before "2022-Q3") ###

$newSheet = $wb.Worksheets.Add($null, $sheetTotal)
$newSheet.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$srcHeaderStyle = $sheetTotal.Cells.Item(1, 2)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i
    $cell = $newSheet.Cells.Item(1, $col)
    $srcHeaderStyle.Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
    $cell.Value = $headers[$i]
}

$srcAStyle = $sheetTotal.Cells.Item(2, 1)
$data = @(
    @{ A = 0; B = "900008"; C = "中信证券稳健回报混合A"; D = "3.38"; E = "86.78"; F = "5.91"; G = "0.1998"; H = 6 },
    @{ A = 1; B = "900078"; C = "中信证券稳健回报混合C"; D = "1.25"; E = "86.78"; F = "5.91"; G = "0.0739"; H = 6 },
    @{ A = 2; B = "005269"; C = "华泰柏瑞港股通量化灵活配置混合"; D = "0.54"; E = "80.96"; F = "1.96"; G = "0.0106"; H = 5 }
)

$r = 2
foreach ($row in $data) {
    $srcAStyle.Copy()
    $newSheet.Cells.Item($r, 1).PasteSpecial(-4122)  # xlPasteFormats
    $newSheet.Cells.Item($r, 1).Value = $row.A

    Set-TextValue $newSheet.Cells.Item($r, 2) $row.B
    $newSheet.Cells.Item($r, 3).Value = $row.C
    Set-TextValue $newSheet.Cells.Item($r, 4) $row.D
    Set-TextValue $newSheet.Cells.Item($r, 5) $row.E
    Set-TextValue $newSheet.Cells.Item($r, 6) $row.F
    Set-TextValue $newSheet.Cells.Item($r, 7) $row.G
    $newSheet.Cells.Item($r, 8).Value = $row.H
    $r++
}

### 3. Match the new sheet's page margins to "总计" (0.75/0.75/1/1/0.5/0.5 in) ###

$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

### 4. Keep "2022-Q3" as the selected/active tab (matches original file) ###

$wb.Worksheets.Item("2022-Q3").Activate()
